# Generate Report for handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet and on each per-locale sheet.
#  - Each per-locale sheet gets its "Latest Target File" / "Latest Handback File"
#    columns filled in (with hyperlinks) and "Latest Handback DateTime" stamped.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# --- Status text: flip every occurrence (Overview + both locale sheets) so the
# shared string is edited in place instead of leaving a stray duplicate. ---
foreach ($ws in $wb.Worksheets) {
    $dims = $ws.UsedRange
    for ($r = 1; $r -le $dims.Rows.Count(); $r++) {
        for ($c = 1; $c -le $dims.Columns.Count(); $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value() -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

# --- zh-cn sheet: rows 2 and 3 now show the handed-back target/handback files ---
$zhTargetUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/886e95aa65cd096e8cd626a5f77b83f4afe0c294/e2e/ba436930-3f19-472c-b819-7d06ea4c6624.md"
$zhTargetName  = "ba436930-3f19-472c-b819-7d06ea4c6624.md"
$zhHandbackUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20ad2656ae378fedefcdccf7188df119e44b8efb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ba436930-3f19-472c-b819-7d06ea4c6624.614134b57d741bace7b01fa8ebab4b2f3c7f6b55.zh-cn.xlf"
$zhHandbackName = "ba436930-3f19-472c-b819-7d06ea4c6624.614134b57d741bace7b01fa8ebab4b2f3c7f6b55.zh-cn.xlf"

$zh.Range("E2").Value = $zhTargetName
$zh.Hyperlinks.Add($zh.Range("E2"), $zhTargetUrl, "", "", $zhTargetName)
$zh.Range("F2").Value = $zhHandbackName
$zh.Hyperlinks.Add($zh.Range("F2"), $zhHandbackUrl, "", "", $zhHandbackName)
$zh.Range("G2").Value = "2016-01-20 08:14:06"

$zh.Range("E3").Value = $zhTargetName
$zh.Hyperlinks.Add($zh.Range("E3"), $zhTargetUrl, "", "", $zhTargetName)
$zh.Range("F3").Value = $zhHandbackName
$zh.Hyperlinks.Add($zh.Range("F3"), $zhHandbackUrl, "", "", $zhHandbackName)
$zh.Range("G3").Value = "2016-01-20 08:14:06"

# --- de-de sheet: same shape, different target locale files/timestamps ---
$deTargetUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/886e95aa65cd096e8cd626a5f77b83f4afe0c294/e2e/ba436930-3f19-472c-b819-7d06ea4c6624.md"
$deTargetName  = "ba436930-3f19-472c-b819-7d06ea4c6624.md"
$deHandbackUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3fca989e1e6b0041375b8dac834195d976f933ad/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ba436930-3f19-472c-b819-7d06ea4c6624.614134b57d741bace7b01fa8ebab4b2f3c7f6b55.de-de.xlf"
$deHandbackName = "ba436930-3f19-472c-b819-7d06ea4c6624.614134b57d741bace7b01fa8ebab4b2f3c7f6b55.de-de.xlf"

$de.Range("E2").Value = $deTargetName
$de.Hyperlinks.Add($de.Range("E2"), $deTargetUrl, "", "", $deTargetName)
$de.Range("F2").Value = $deHandbackName
$de.Hyperlinks.Add($de.Range("F2"), $deHandbackUrl, "", "", $deHandbackName)
$de.Range("G2").Value = "2016-01-20 08:14:27"

$de.Range("E3").Value = $deTargetName
$de.Hyperlinks.Add($de.Range("E3"), $deTargetUrl, "", "", $deTargetName)
$de.Range("F3").Value = $deHandbackName
$de.Hyperlinks.Add($de.Range("F3"), $deHandbackUrl, "", "", $deHandbackName)
$de.Range("G3").Value = "2016-01-20 08:14:27"

# --- Match the hyperlink-column look (underlined custom HyperLink style) used
# by columns A/C so the newly-populated E/F columns are styled the same way. ---
$zh.Range("E2:F3").Style = "HyperLink"
$de.Range("E2:F3").Style = "HyperLink"

Write-Host "Handback report generated."
